$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 552.125
$ws.Range("I4").Value = 538.93335
$ws.Range("J4").Value = 750
$ws.Range("K4").Value = 538.93335
$ws.Range("L4").Value = 750
$ws.Range("M4").Value = -424.93335
$ws.Range("N4").Value = -978
$ws.Range("H15").Value = 533.1091
$ws.Range("I15").Value = 533.1091
$ws.Range("K15").Value = 1599.3273
$ws.Range("M15").Value = -1430.3273
$ws.Range("H98").Value = 2058.476
$ws.Range("J98").Value = 1605
$ws.Range("L98").Value = 1605
$ws.Range("N98").Value = -4601
$ws.Range("H112").Value = 201416.6
$ws.Range("J112").Value = 201416.6
$ws.Range("L112").Value = 604249.8
$ws.Range("N112").Value = -606465.8
$ws.Range("H116").Value = 9799.786
$ws.Range("J116").Value = 9784.385
$ws.Range("L116").Value = 9784.385
$ws.Range("N116").Value = -16668.385
$ws.Range("H122").Value = 2058.476
$ws.Range("J122").Value = 1605
$ws.Range("L122").Value = 4815
$ws.Range("N122").Value = -9715
$ws.Range("H129").Value = 1687.2307
$ws.Range("I129").Value = 1086.5
$ws.Range("J129").Value = 1954.2222
$ws.Range("K129").Value = 3259.5
$ws.Range("L129").Value = 5862.6666
$ws.Range("M129").Value = 1740.5
$ws.Range("N129").Value = -15862.6666
$ws.Range("H138").Value = 6120.064
$ws.Range("I138").Value = 11765.706
$ws.Range("J138").Value = 2920.8667
$ws.Range("K138").Value = 35297.118
$ws.Range("L138").Value = 8762.6001
$ws.Range("M138").Value = -30157.118
$ws.Range("N138").Value = -19042.6001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 170959.2
$ws.Range("I32").Value = 187086.9
$ws.Range("J32").Value = 25809.834
$ws.Range("K32").Value = 187086.9
$ws.Range("L32").Value = 25809.834
$ws.Range("M32").Value = -186799.9
$ws.Range("N32").Value = -26383.834
$ws.Range("H74").Value = 2465.365
$ws.Range("I74").Value = 1461.909
$ws.Range("J74").Value = 3569.1667
$ws.Range("K74").Value = 1461.909
$ws.Range("L74").Value = 3569.1667
$ws.Range("M74").Value = -587.9090000000001
$ws.Range("N74").Value = -5317.1667
$ws.Range("H77").Value = 2465.365
$ws.Range("I77").Value = 1461.909
$ws.Range("J77").Value = 3569.1667
$ws.Range("K77").Value = 7309.545
$ws.Range("L77").Value = 17845.8335
$ws.Range("M77").Value = -2941.545
$ws.Range("N77").Value = -26581.8335
$ws.Range("H102").Value = 2529.75
$ws.Range("I102").Value = 2494.1667
$ws.Range("K102").Value = 2494.1667
$ws.Range("M102").Value = -872.1667000000002
$ws.Range("H132").Value = 2522.3333
$ws.Range("I132").Value = 1451.9166
$ws.Range("K132").Value = 4355.7498
$ws.Range("M132").Value = -1825.7498

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5955.8823
$ws.Range("I86").Value = 2372.125
$ws.Range("J86").Value = 9141.444
$ws.Range("K86").Value = 2372.125
$ws.Range("L86").Value = 9141.444
$ws.Range("M86").Value = -1249.125
$ws.Range("N86").Value = -11387.444
$ws.Range("H89").Value = 5955.8823
$ws.Range("I89").Value = 2372.125
$ws.Range("J89").Value = 9141.444
$ws.Range("K89").Value = 11860.625
$ws.Range("L89").Value = 45707.22
$ws.Range("M89").Value = -6244.625
$ws.Range("N89").Value = -56939.22
$ws.Range("H105").Value = 7134.24
$ws.Range("I105").Value = 8303.5625
$ws.Range("K105").Value = 8303.5625
$ws.Range("M105").Value = -6556.5625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5006933.5
$ws.Range("H31").Value = 2950.6135
$ws.Range("J31").Value = 3043.8462
$ws.Range("L31").Value = 3043.8462
$ws.Range("N31").Value = -3633.8462
$ws.Range("H34").Value = 2950.6135
$ws.Range("J34").Value = 3043.8462
$ws.Range("L34").Value = 3043.8462
$ws.Range("N34").Value = -3447.8462
$ws.Range("H58").Value = 2869.25
$ws.Range("I58").Value = 1770.5714
$ws.Range("K58").Value = 1770.5714
$ws.Range("M58").Value = -1567.5714
$ws.Range("H107").Value = 1450.1562
$ws.Range("I107").Value = 1548.75
$ws.Range("J107").Value = 1351.5625
$ws.Range("K107").Value = 1548.75
$ws.Range("L107").Value = 1351.5625
$ws.Range("M107").Value = 371.25
$ws.Range("N107").Value = -5191.5625
$ws.Range("H132").Value = 17547610
$ws.Range("I132").Value = 4142.1
$ws.Range("K132").Value = 12426.3
$ws.Range("M132").Value = -9896.300000000001
$ws.Range("H134").Value = 2363
$ws.Range("I134").Value = 1875.8
$ws.Range("K134").Value = 5627.4
$ws.Range("M134").Value = -3092.4
$ws.Range("H136").Value = 2869.25
$ws.Range("I136").Value = 1770.5714
$ws.Range("K136").Value = 5311.7142
$ws.Range("M136").Value = -2761.7142
$ws.Range("H133").Value = 29326
$ws.Range("J133").Value = 29326
$ws.Range("L133").Value = 29326
$ws.Range("N133").Value = -34386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3782887.5
$ws.Range("I4").Value = 4526416.5
$ws.Range("K4").Value = 13579249.5
$ws.Range("M4").Value = -13579137.5
$ws.Range("H34").Value = 590.06665
$ws.Range("J34").Value = 1022
$ws.Range("L34").Value = 3066
$ws.Range("N34").Value = -3234
$ws.Range("H107").Value = 40000384
$ws.Range("J107").Value = 58823828
$ws.Range("L107").Value = 176471484
$ws.Range("N107").Value = -176475324
$ws.Range("H140").Value = 22730048
$ws.Range("I140").Value = 31252400
$ws.Range("J140").Value = 3777.6667
$ws.Range("K140").Value = 93757200
$ws.Range("L140").Value = 11333.0001
$ws.Range("M140").Value = -93752020
$ws.Range("N140").Value = -21693.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 791.2381
$ws.Range("I97").Value = 729.7692
$ws.Range("J97").Value = 891.125
$ws.Range("K97").Value = 729.7692
$ws.Range("L97").Value = 891.125
$ws.Range("M97").Value = -233.7692
$ws.Range("N97").Value = -1883.125
$ws.Range("H132").Value = 13128657
$ws.Range("I132").Value = 1556
$ws.Range("K132").Value = 4668
$ws.Range("M132").Value = -2138

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4498.933
$ws.Range("J22").Value = 6144.7
$ws.Range("L22").Value = 6144.7
$ws.Range("N22").Value = -6734.7
$ws.Range("H27").Value = 4498.933
$ws.Range("J27").Value = 6144.7
$ws.Range("L27").Value = 6144.7
$ws.Range("N27").Value = -6358.7
$ws.Range("H55").Value = 533.73334
$ws.Range("J55").Value = 651.1818
$ws.Range("L55").Value = 651.1818
$ws.Range("N55").Value = -997.1818
$ws.Range("H61").Value = 2817.639
$ws.Range("I61").Value = 2682.6553
$ws.Range("K61").Value = 2682.6553
$ws.Range("M61").Value = -2480.6553
$ws.Range("H82").Value = 4649
$ws.Range("I82").Value = 1054.375
$ws.Range("K82").Value = 1054.375
$ws.Range("M82").Value = -693.375
$ws.Range("H85").Value = 4649
$ws.Range("I85").Value = 1054.375
$ws.Range("K85").Value = 1054.375
$ws.Range("M85").Value = 193.625
$ws.Range("H113").Value = 2817.639
$ws.Range("I113").Value = 2682.6553
$ws.Range("K113").Value = 2682.6553
$ws.Range("M113").Value = -512.6552999999999
$ws.Range("H122").Value = 4635.273
$ws.Range("I122").Value = 2921.6924
$ws.Range("K122").Value = 8765.0772
$ws.Range("M122").Value = -6315.0772
$ws.Range("H132").Value = 11798.533
$ws.Range("J132").Value = 14634.363
$ws.Range("L132").Value = 43903.089
$ws.Range("N132").Value = -48963.089

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 9999999
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H45").Value = 12880.4
$ws.Range("I45").Value = 6495.6665
$ws.Range("J45").Value = 15616.714
$ws.Range("K45").Value = 6495.6665
$ws.Range("L45").Value = 15616.714
$ws.Range("M45").Value = -6004.6665
$ws.Range("N45").Value = -16598.714
$ws.Range("H113").Value = 731.3333
$ws.Range("I113").Value = 499.33334
$ws.Range("J113").Value = 847.3333
$ws.Range("K113").Value = 1498.00002
$ws.Range("L113").Value = 2541.9999
$ws.Range("M113").Value = 671.9999800000001
$ws.Range("N113").Value = -6881.9999
$ws.Range("H122").Value = 1555.9678
$ws.Range("J122").Value = 599.75
$ws.Range("L122").Value = 1799.25
$ws.Range("N122").Value = -6699.25
$ws.Range("H132").Value = 27679.658
$ws.Range("I132").Value = 35554.55
$ws.Range("K132").Value = 106663.65
$ws.Range("M132").Value = -104133.65
$ws.Range("H136").Value = 18572.932
$ws.Range("I136").Value = 28734.139
$ws.Range("J136").Value = 2668.4348
$ws.Range("K136").Value = 86202.417
$ws.Range("L136").Value = 8005.3044
$ws.Range("M136").Value = -83652.417
$ws.Range("N136").Value = -13105.3044
